# Update column F (dSF) values on Sheet1 per repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -7
    5  = -2
    6  = 0
    12 = -2
    15 = 0
    20 = -1
    22 = -1
    27 = 2
    31 = -1
    33 = 1
    35 = 1
    36 = 0
    39 = -3
    41 = 2
    43 = 1
    47 = 4
    51 = 1
    57 = 1
    70 = 2
    75 = -4
    77 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
